$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1779.6428
$ws.Range("I9").Value = 2032.8572
$ws.Range("J9").Value = 1526.4286
$ws.Range("K9").Value = 2032.8572
$ws.Range("L9").Value = 1526.4286
$ws.Range("M9").Value = -1863.8572
$ws.Range("N9").Value = -1864.4286

$ws.Range("H20").Value = 1746.2858
$ws.Range("I20").Value = 700
$ws.Range("J20").Value = 8024
$ws.Range("K20").Value = 700
$ws.Range("L20").Value = 8024
$ws.Range("M20").Value = -470
$ws.Range("N20").Value = -8484

$ws.Range("H35").Value = 1746.2858
$ws.Range("I35").Value = 700
$ws.Range("J35").Value = 8024
$ws.Range("K35").Value = 700
$ws.Range("L35").Value = 8024
$ws.Range("M35").Value = -321
$ws.Range("N35").Value = -8782

$ws.Range("H40").Value = 4333.3335
$ws.Range("I40").Value = 4333.3335
$ws.Range("K40").Value = 4333.3335
$ws.Range("M40").Value = -4158.3335

$ws.Range("H100").Value = 3721.5557
$ws.Range("I100").Value = 2499.6667
$ws.Range("J100").Value = 4332.5
$ws.Range("K100").Value = 2499.6667
$ws.Range("L100").Value = 4332.5
$ws.Range("M100").Value = -1958.6667
$ws.Range("N100").Value = -5414.5

$ws.Range("H107").Value = 608.63635
$ws.Range("I107").Value = 509
$ws.Range("K107").Value = 509
$ws.Range("M107").Value = 1411

$ws.Range("H113").Value = 10694.667
$ws.Range("J113").Value = 10328.571
$ws.Range("L113").Value = 10328.571
$ws.Range("N113").Value = -16836.571

$ws.Range("H129").Value = 765.6667
$ws.Range("I129").Value = 765.6667
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2297.0001
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2702.9999
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 100398
$ws.Range("J114").Value = 100398
$ws.Range("L114").Value = 100398
$ws.Range("N114").Value = -109076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1187.5
$ws.Range("I22").Value = 1187.5
$ws.Range("K22").Value = 1187.5
$ws.Range("M22").Value = -1014.5

$ws.Range("H94").Value = 2481.0156
$ws.Range("J94").Value = 3148.2666
$ws.Range("L94").Value = 3148.2666
$ws.Range("N94").Value = -4050.2666

$ws.Range("H105").Value = 1156.2632
$ws.Range("I105").Value = 810.625
$ws.Range("K105").Value = 810.625
$ws.Range("M105").Value = 936.375

$ws.Range("H134").Value = 9263.6
$ws.Range("I134").Value = 1850.0869
$ws.Range("K134").Value = 5550.2607
$ws.Range("M134").Value = -3015.2607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 163.78947
$ws.Range("I7").Value = 93.916664
$ws.Range("J7").Value = 283.57144
$ws.Range("K7").Value = 93.916664
$ws.Range("L7").Value = 283.57144
$ws.Range("M7").Value = 19.083336
$ws.Range("N7").Value = -509.57144

$ws.Range("H39").Value = 6525
$ws.Range("I39").Value = 6525
$ws.Range("K39").Value = 6525
$ws.Range("M39").Value = -6134

$ws.Range("H49").Value = 6525
$ws.Range("I49").Value = 6525
$ws.Range("K49").Value = 6525
$ws.Range("M49").Value = -6343

$ws.Range("H118").Value = 47500
$ws.Range("J118").Value = 47500
$ws.Range("L118").Value = 47500
$ws.Range("N118").Value = -50814

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1182.8462
$ws.Range("I11").Value = 1505.375
$ws.Range("J11").Value = 666.8
$ws.Range("K11").Value = 4516.125
$ws.Range("L11").Value = 2000.4
$ws.Range("M11").Value = -4376.125
$ws.Range("N11").Value = -2280.4

$ws.Range("H13").Value = 333.9
$ws.Range("I13").Value = 93.333336
$ws.Range("J13").Value = 694.75
$ws.Range("K13").Value = 280.000008
$ws.Range("L13").Value = 2084.25
$ws.Range("M13").Value = -112.000008
$ws.Range("N13").Value = -2420.25

$ws.Range("H87").Value = 10250
$ws.Range("I87").Value = 4000
$ws.Range("J87").Value = 16500
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 49500
$ws.Range("M87").Value = -10752
$ws.Range("N87").Value = -51996

$ws.Range("H90").Value = 10250
$ws.Range("I90").Value = 4000
$ws.Range("J90").Value = 16500
$ws.Range("K90").Value = 36000
$ws.Range("L90").Value = 148500
$ws.Range("M90").Value = -29760
$ws.Range("N90").Value = -160980

$ws.Range("H114").Value = 1179.75
$ws.Range("J114").Value = 999
$ws.Range("L114").Value = 2997
$ws.Range("N114").Value = -9505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 299987
$ws.Range("J105").Value = 299987
$ws.Range("L105").Value = 299987
$ws.Range("N105").Value = -306975

$ws.Range("H107").Value = 655.0454999999999
$ws.Range("I107").Value = 235.625
$ws.Range("K107").Value = 235.625
$ws.Range("M107").Value = 1684.375

$ws.Range("H111").Value = 39146.5
$ws.Range("J111").Value = 39146.5
$ws.Range("L111").Value = 39146.5
$ws.Range("N111").Value = -45280.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 8667
$ws.Range("I13").Value = 12750.5
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 12750.5
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = -12610.5
$ws.Range("N13").Value = -780

$ws.Range("H25").Value = 978498.75
$ws.Range("J25").Value = 1303998.4
$ws.Range("L25").Value = 1303998.4
$ws.Range("N25").Value = -1304458.4

$ws.Range("H42").Value = 24105.166
$ws.Range("J42").Value = 24105.166
$ws.Range("L42").Value = 24105.166
$ws.Range("N42").Value = -25231.166

$ws.Range("H49").Value = 24105.166
$ws.Range("J49").Value = 24105.166
$ws.Range("L49").Value = 24105.166
$ws.Range("N49").Value = -24399.166

$ws.Range("H68").Value = 3546.3333
$ws.Range("J68").Value = 3950
$ws.Range("L68").Value = 3950
$ws.Range("N68").Value = -5448

$ws.Range("H71").Value = 3546.3333
$ws.Range("J71").Value = 3950
$ws.Range("L71").Value = 19750
$ws.Range("N71").Value = -27238

$ws.Range("H93").Value = 10901.579
$ws.Range("I93").Value = 13031.667
$ws.Range("J93").Value = 7250
$ws.Range("K93").Value = 13031.667
$ws.Range("L93").Value = 7250
$ws.Range("M93").Value = -11783.667
$ws.Range("N93").Value = -9746

$ws.Range("H100").Value = 3615.9565
$ws.Range("I100").Value = 3451.0588
$ws.Range("K100").Value = 3451.0588
$ws.Range("M100").Value = -2910.0588

$ws.Range("H137").Value = 84449.45
$ws.Range("J137").Value = 84525.734
$ws.Range("L137").Value = 84525.734
$ws.Range("N137").Value = -94725.734

$ws.Range("H139").Value = 58350.91
$ws.Range("J139").Value = 58350.91
$ws.Range("L139").Value = 58350.91
$ws.Range("N139").Value = -68630.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2020.7222
$ws.Range("I96").Value = 1591
$ws.Range("K96").Value = 1591
$ws.Range("M96").Value = -218

$ws.Range("H132").Value = 5509.183
$ws.Range("I132").Value = 2989.0952
$ws.Range("K132").Value = 8967.285600000001
$ws.Range("M132").Value = -6437.285600000001

$ws.Range("H141").Value = 68166
$ws.Range("J141").Value = 68166
$ws.Range("L141").Value = 68166
$ws.Range("N141").Value = -78526
